$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (reordered/extended) monthly data for columns A-D, rows 2-68.
# Source row/commit reorders each year so Oct/Nov/Dec precede Jan-Sep, and
# extends the series from 2021-12 through 2023-07 (new years 2022 & 2023).
$data = @(
  @("2018-10", 101.6, 107.9, 99.90000000000001),
  @("2018-11", 100.9, 106.6, 98.7),
  @("2018-12", 100.8, 106.6, 99.2),
  @("2018-01", 105.7415, 116.889, 110.4718),
  @("2018-02", 105, 114.2, 109.6),
  @("2018-03", 105.9, 113.2, 111.4),
  @("2018-04", 106.2, 111.5, 113.4),
  @("2018-05", 106.9, 110.7, 115.2),
  @("2018-06", 106.7, 112.9, 114.4),
  @("2018-07", 106.3, 113.8, 112.7),
  @("2018-08", 105.6, 115.1, 109.7),
  @("2018-09", 103.7, 111.7, 104.7),
  @("2019-10", 97.7, 82.09999999999999, 90.5),
  @("2019-11", 97.8, 80.5, 92.5),
  @("2019-12", 98, 79.5, 94.2),
  @("2019-01", 100.3, 101.9, 98.5),
  @("2019-02", 100.1, 99, 97.59999999999999),
  @("2019-03", 99.2, 99.7, 95.5),
  @("2019-04", 99.2, 101, 93.90000000000001),
  @("2019-05", 98.7, 100.3, 91.7),
  @("2019-06", 98.09999999999999, 95.3, 90),
  @("2019-07", 98.09999999999999, 89, 89.3),
  @("2019-08", 97.7, 85.5, 88.7),
  @("2019-09", 97.40000000000001, 82.59999999999999, 88.8),
  @("2020-10", 99.09999999999999, 93.59999999999999, 97.8),
  @("2020-11", 99.5, 96, 98.3),
  @("2020-12", 99.5, 96.8, 97.8),
  @("2020-01", 98, 81.90000000000001, 95.7),
  @("2020-02", 98.3, 83.5, 96.40000000000001),
  @("2020-03", 98.90000000000001, 84.7, 98.09999999999999),
  @("2020-04", 98.2, 84, 95.40000000000001),
  @("2020-05", 97.3, 84, 93),
  @("2020-06", 97.7, 86.7, 93.2),
  @("2020-07", 97.90000000000001, 92.2, 95.59999999999999),
  @("2020-08", 98.5, 94.59999999999999, 97.40000000000001),
  @("2020-09", 98.90000000000001, 95.59999999999999, 97.59999999999999),
  @("2021-10", 102.9, 126.1, 107.7),
  @("2021-11", 103.1, 124.3, 108.9),
  @("2021-12", 103.6, 121, 108),
  @("2021-01", 99.7, 101.4, 99.40000000000001),
  @("2021-02", 100.2, 106.8, 101.5),
  @("2021-03", 100.9, 110.7, 103.2),
  @("2021-04", 102.1, 112.6, 107.1),
  @("2021-05", 103.2, 119.6, 110.4),
  @("2021-06", 103.2, 123.6, 111.1),
  @("2021-07", 103.5, 120.7, 108.1),
  @("2021-08", 103.3, 124.6, 107.1),
  @("2021-09", 103.1, 123.9, 107.1),
  @("2022-10", 101.1, 121.1, 97.8),
  @("2022-11", 100.9, 124.2, 95.8),
  @("2022-12", 100.3, 123.5, 95.90000000000001),
  @("2022-01", 103.2, 121.5, 105.8),
  @("2022-02", 102.5, 116.8, 103.3),
  @("2022-03", 101.4, 113.7, 100),
  @("2022-04", 101.2, 111.9, 100.3),
  @("2022-05", 100.9, 111.9, 100),
  @("2022-06", 101.4, 114.5, 99.59999999999999),
  @("2022-07", 100.8, 120.3, 100.4),
  @("2022-08", 101.1, 117.9, 100),
  @("2022-09", 101.2, 119.7, 98.7),
  @("2023-01", 100.3, 120.7, 96.40000000000001),
  @("2023-02", 100.3, 119.2, 95.8),
  @("2023-03", 100, 111.8, 94.40000000000001),
  @("2023-04", 99.2, 103.9, 92.40000000000001),
  @("2023-05", 98.40000000000001, 92.3, 91.40000000000001),
  @("2023-06", 97.7, 81.40000000000001, 90.09999999999999),
  @("2023-07", 97.8, 76.2, 89)
)

# Column A (the date label) already carries a bold/centered/bordered style
# (style index 1) on every existing data row; reuse it as the paste-format
# source so newly appended rows (50-68) match the existing look.
$styleSrc = $ws.Range("A2")
$existingMaxRow = 49

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  if ($r -gt $existingMaxRow) {
    $styleSrc.Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
  }
  $r = $r + 1
}
